$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.602.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4668"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9020"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.853.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.252"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.327"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008564"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.639.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.007"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.936"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("E28").Value = "  -2.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "

$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08761"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.141"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.765"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7353"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.431"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.076"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01933"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.40%  "

$ws.Range("E39").Value = "  +1.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05110"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5076"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.802"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.990"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4686"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.566"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
